# Generate Report for Handback
#
# This models a handback run where the zh-cn and de-de handback transforms
# failed because the handback file name did not match the handoff file name.
# Every cell that previously read "Ready for handoff" for the
# 07d6f1e9-fc20-4950-8cf1-c67ba3aee511 row now reads "Handback transform
# failed" (Overview's per-locale status columns plus each locale sheet's own
# Status column), and the per-locale sheets additionally get an Error Detail
# value recorded for that row (with the Error Detail column widened so the
# message is readable).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the 07d6f1e9... row as failed instead of "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: same status flip, plus widen the Error Detail column (P) and record the error ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsZhCn.Range("P3").Value = "Handback file name: lo4nnokq.n15 is different with handoff file name: 07d6f1e9-fc20-4950-8cf1-c67ba3aee511.69d0bb56069892976b6632d4e2837d04369e0d72.zh-cn."

# --- de-de sheet: same status flip, plus widen the Error Detail column (P) and record the error ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Range("P3").Value = "Handback file name: lo4nnokq.n15 is different with handoff file name: 07d6f1e9-fc20-4950-8cf1-c67ba3aee511.69d0bb56069892976b6632d4e2837d04369e0d72.de-de."
